# The commit removes the slide "Tips by Payment Method" (old slide index 7,
# p:sldId 270) from the deck. All other slides keep their original content
# and relative order; the sldIdLst entries after it simply shift up by one
# position (which Slides.Item(7).Delete() reproduces natively).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$s.Delete()
